# Update the "想去人数" (want-to-go count) figures in column F for the
# exhibition listings on both the "展览" and "全部类型" worksheets.
# These are small incremental refreshes pulled from the live data source
# (bilibili event pages), matching the gh-pages regeneration at 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        3  = 27
        4  = 1387
        5  = 316
        7  = 10704
        8  = 23
        10 = 294
        11 = 1037
        12 = 711
        13 = 12052
        14 = 12503
    }
    "全部类型" = @{
        4  = 27
        5  = 1387
        6  = 316
        8  = 10704
        9  = 23
        11 = 294
        12 = 1037
        13 = 711
        14 = 12052
        15 = 12503
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowMap[$row]
    }
}
